# Edit TransectsDistanceTutorial.xlsx:
#  - Shorten header labels in row 1 (A1, B1, C1)
#  - Move the active cell / selection on the frozen pane to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text for the three columns whose shared-string
# values changed in the diff.
$ws.Range("A1").Value = "Primary observer"
$ws.Range("B1").Value = "Secondary observer"
$ws.Range("C1").Value = "Transect length "

# Update the selection / active cell for the bottom-left (frozen) pane
# from A2 to C2.
$ws.Range("C2").Select()
